# Auto-generated Excel COM-interop edit script
# Refreshes the scraped crypto price/volume table (and fixes the row order for two
# coins whose ranking swapped) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '64.412.51'
$ws.Range('E2').Value = '  -2.45%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.179.92'
$ws.Range('E3').Value = '  -3.98%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.46%  '

# Row 6: Solana
$ws.Range('E6').Value = '  -7.45%  '

# Row 7: XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.609'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.75%  '

# Row 8: USDC
$ws.Range('E8').Value = '  -0.13%  '

# Row 9: LidoStakedEther
$ws.Range('D9').Value = '3.188.68'
$ws.Range('E9').Value = '  -3.56%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  -3.58%  '

# Row 11: Toncoin
$ws.Range('E11').Value = '  -0.27%  '

# Row 12: Cardano
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.388'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.09%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range('D13').Value = '3.730.30'
$ws.Range('E13').Value = '  -4.01%  '

# Row 14: TRON
$ws.Range('E14').Value = '  -1.61%  '

# Row 15: WrappedBTC
$ws.Range('D15').Value = '64.472.08'
$ws.Range('E15').Value = '  -2.42%  '

# Row 16: Avalanche
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.31'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.09%  '

# Row 17: ShibaInu
$ws.Range('E17').Value = '  -3.40%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '3.183.38'
$ws.Range('E18').Value = '  -3.85%  '

# Row 19: BitcoinCash
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '418.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.28%  '

# Row 20: Chainlink
$ws.Range('E20').Value = '  -1.20%  '

# Row 21: Polkadot
$ws.Range('E21').Value = '  -3.13%  '

# Row 22: Uniswap
$ws.Range('E22').Value = '  -3.44%  '

# Row 23: Dai
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.15%  '

# Row 24: Litecoin
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.34'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.91%  '

# Row 25: LEO
$ws.Range('E25').Value = '  +0.02%  '

# Row 26: Kaspa
$ws.Range('E26').Value = '  +2.72%  '

# Row 27: Polygon
$ws.Range('E27').Value = '  -4.28%  '

# Row 28: PEPE
$ws.Range('E28').Value = '  -6.25%  '

# Row 29: InternetComputer(DFINITY)
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.75'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.49%  '

# Row 30: Binance-PegBSC-USD
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.995'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.86%  '

# Row 31: PancakeSwap
$ws.Range('E31').Value = '  -3.38%  '

# Row 32: EthereumClassic
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.74'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.66%  '

# Row 33: USDe
$ws.Range('E33').Value = '  -0.10%  '

# Row 34: NEARProtocol
$ws.Range('E34').Value = '  -1.92%  '

# Row 35: Aptos
$ws.Range('E35').Value = '  -2.82%  '

# Row 36: Fetch.AI
$ws.Range('E36').Value = '  -3.40%  '

# Row 37: Monero
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '156.94'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.19%  '

# Row 38: ImmutableX
$ws.Range('E38').Value = '  -4.69%  '

# Row 39: Stacks
$ws.Range('E39').Value = '  -4.88%  '

# Row 40: Maker
$ws.Range('D40').Value = '2.697.30'
$ws.Range('E40').Value = '  -6.02%  '

# Row 41: Filecoin
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.22'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.69%  '

# Row 42: EnergySwap
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '24.22'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.64%  '

# Row 43: OKB
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.23'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.62%  '

# Row 44: Mantle
$ws.Range('E44').Value = '  -5.41%  '

# Row 45: Hedera
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0622'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.56%  '

# Row 46: RenderToken
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.57'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.55%  '

# Row 47: VeChain
$ws.Range('E47').Value = '  -2.41%  '

# Row 48: Bittensor
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '292.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.71%  '

# Row 49: InjectiveProtocol
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '21.40'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.10%  '

# Row 50: Stellar
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0992'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.32%  '

# Row 51: FirstDigitalUSD
$ws.Range('E51').Value = '  -0.19%  '
